$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Target change: the paragraph "(Processos extraídos das entrevistas)"
# gets split into two runs - "(" (unchanged formatting) and
# "Processos extraídos das entrevistas)" (font size bumped to 12pt) - and
# the existing "_GoBack" bookmark (previously sitting mid-sentence in the
# "Compra de veículo" paragraph, right after "realizado ") is moved to the
# end of this paragraph, right after the closing ")".
# -----------------------------------------------------------------------

# --- Step 1: locate the paragraph "(Processos extraídos das entrevistas)" ---
$findRange = $d.Content
$found = $findRange.Find.Execute(
    "(Processos extraídos das entrevistas)",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

if ($found) {
    $openParenStart = $findRange.Start   # position right before "("
    $closeParenEnd  = $findRange.End     # position right after ")"

    # --- Step 2: split the single run into "(" + "Processos extraídos das
    #     entrevistas)" by raising the font size of everything after "(" ---
    $secondPartRange = $d.Range($openParenStart + 1, $closeParenEnd)
    $secondPartRange.Font.Size = 12

    # --- Step 3: remove the old "_GoBack" bookmark, wherever it is now ---
    if ($d.Bookmarks.Exists("_GoBack")) {
        $oldBookmark = $d.Bookmarks("_GoBack")
        $oldBookmark.Delete()
    }

    # --- Step 4: re-create "_GoBack" right after the ")" that now ends
    #     this paragraph.
    #     Note: adding a bookmark on a collapsed Range that sits exactly at
    #     a paragraph's last character offset (immediately before the
    #     paragraph mark) snaps to the wrong location, so a placeholder
    #     character is inserted first, the bookmark is added next to it,
    #     and the placeholder is removed afterwards. ---
    $placeholderRange = $d.Range($closeParenEnd, $closeParenEnd)
    $placeholderRange.InsertAfter("X")

    $bookmarkRange = $d.Range($closeParenEnd, $closeParenEnd)
    $d.Bookmarks.Add("_GoBack", $bookmarkRange)

    $placeholderCleanupRange = $d.Range($closeParenEnd, $closeParenEnd + 1)
    $placeholderCleanupRange.Text = ""
}
